$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header volume number (35 -> 36) ---
$ws.Range("A8").Value = "Volume 30   Number  36"

# --- Update report week covering dates ---
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Simple value updates (style already correct) ---
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = -7.692307692307
$ws.Range("L15").Value = -7.692307692307
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -29.411764705882
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 98
$ws.Range("K16").Value = 25.641025641025
$ws.Range("L16").Value = 36.111111111111
$ws.Range("M16").Value = -20.967741935483
$ws.Range("N16").Value = -81.918819188191
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 38.461538461538
$ws.Range("I17").Value = 152
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 26.666666666666
$ws.Range("L17").Value = 20.63492063492
$ws.Range("M17").Value = 92.405063291139
$ws.Range("N17").Value = -27.619047619047
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 128
$ws.Range("K18").Value = 21.09375
$ws.Range("L18").Value = 63.157894736842
$ws.Range("M18").Value = -21.319796954314
$ws.Range("N18").Value = -85.432330827067
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = -52
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 81
$ws.Range("H19").Value = -53.086419753086
$ws.Range("I19").Value = 459
$ws.Range("J19").Value = 515
$ws.Range("K19").Value = -10.873786407767
$ws.Range("L19").Value = 17.994858611825
$ws.Range("M19").Value = 59.375
$ws.Range("N19").Value = -13.068181818181
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 80
$ws.Range("L20").Value = 82.258064516129
$ws.Range("M20").Value = -3.418803418803
$ws.Range("N20").Value = -90.93825180433
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -21.705426356589
$ws.Range("I21").Value = 991
$ws.Range("J21").Value = 958
$ws.Range("K21").Value = 3.444676409185
$ws.Range("L21").Value = 30.738786279683
$ws.Range("M21").Value = 21.001221001221
$ws.Range("N21").Value = -72.586445366528
$ws.Range("L22").Value = 12.5
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -27.5
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 162
$ws.Range("H24").Value = -12.345679012345
$ws.Range("I24").Value = 1261
$ws.Range("J24").Value = 1265
$ws.Range("K24").Value = -0.316205533596
$ws.Range("L24").Value = 51.5625
$ws.Range("M24").Value = 86.814814814814
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -43.75
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 10.81081081081
$ws.Range("I25").Value = 403
$ws.Range("J25").Value = 315
$ws.Range("K25").Value = 27.936507936507
$ws.Range("L25").Value = 36.148648648648
$ws.Range("M25").Value = 31.270358306188
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 6.666666666666
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = 16.666666666666

# --- Cells needing value type/style change (restyle) ---
# Step 1: set values (text values prefixed with apostrophe to force text type)
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "'***.*"
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C26").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "'***.*"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "'***.*"

# Step 2: copy correct number format/style onto those cells from safe reference cells
$ws.Range("C36").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C34").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("K36").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
